# Vendor_Pincode_Mapping_Template rework:
# The sheet used to hold 10 vendor-upload columns (Vendor Name, Vendor_ID,
# Appliance, Appliance_ID, Brand, Area, Pincode, Region, City, State) with
# a placeholder row underneath. The new template only needs two columns,
# "Pincode" then "Appliance", bold+centered headers and centered sample
# placeholders underneath.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the now-unused columns C:K (Appliance_ID, Brand, Area, Region, City,
# State, ...) and shift everything left.
$ws.Range("C1:K1").EntireColumn.Delete()

# Column A becomes Pincode, column B becomes Appliance.
$ws.Range("A1").Value = "Pincode"
$ws.Range("B1").Value = "Appliance"
$ws.Range("A2").Value = "{vendor:Pincode}"
$ws.Range("B2").Value = "{vendor:Appliance}"

# Header row: bold + centered.
$headerRange = $ws.Range("A1:B1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter

# Placeholder row: centered, regular weight.
$dataRange = $ws.Range("A2:B2")
$dataRange.Font.Bold = $false
$dataRange.HorizontalAlignment = -4108  # xlCenter

# Widen the two remaining columns to fit the new content.
$ws.Columns.Item(1).ColumnWidth = 27.25
$ws.Columns.Item(2).ColumnWidth = 98.75

# Put the selection back on A1 (D13 no longer exists on the trimmed sheet).
[void]$ws.Range("A1").Select()
